# Insert a new weekly record row for "Apio" (Terminal Hortofrutícola Agro
# Chillán) at row 92 of the data table, pushing the existing rows 92-164
# down to 93-165 (row 164 -> row 165, extending the used range to A1:R165).
# The new row captures a later date (2022-01-05, serial 44566) than any
# existing record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 92, shifting rows 92:164 down
# to 93:165 (Excel copies the formatting of the row above, matching the
# date-style cell already present on the neighbouring rows' column D).
$ws.Rows(92).Insert()

$ws.Range("A92").Value = 7
$ws.Range("B92").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C92").Value = 'Ñuble'
$ws.Range("D92").Value = 44566
$ws.Range("E92").Value = 16
$ws.Range("F92").Value = 100112017
$ws.Range("G92").Value = 'Apio'
$ws.Range("H92").Value = 'Americana (o)'
$ws.Range("I92").Value = 'Primera'
$ws.Range("J92").Value = 60
$ws.Range("K92").Value = 8000
$ws.Range("L92").Value = 8500
$ws.Range("M92").Value = 8250
$ws.Range("N92").Value = '$/docena de matas'
$ws.Range("O92").Value = 'Provincia del Elquí'
$ws.Range("P92").Value = 1375
$ws.Range("Q92").Value = 6
$ws.Range("R92").Value = 'Hortaliza'
